$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.477.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("E2").Style = "Normal"

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.610.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +10.27%  "
$ws.Range("E3").Style = "Normal"

# Row 4 - TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("E5").Style = "Normal"

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("E6").Style = "Normal"

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.604"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +6.07%  "
$ws.Range("E7").Style = "Normal"

# Row 8 - USDC
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E8").Style = "Normal"

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +12.41%  "
$ws.Range("E9").Style = "Normal"

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +12.63%  "
$ws.Range("E10").Style = "Normal"

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.79%  "
$ws.Range("E11").Style = "Normal"

# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +15.59%  "
$ws.Range("E12").Style = "Normal"

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.011.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +10.45%  "
$ws.Range("E13").Style = "Normal"

# Row 14 - TRON
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("E14").Style = "Normal"

# Row 15 - WrappedEther
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.608.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +10.38%  "
$ws.Range("E15").Style = "Normal"

# Row 16 - Polygon
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.906"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +11.29%  "
$ws.Range("E16").Style = "Normal"

# Row 17 - Chainlink
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +9.34%  "
$ws.Range("E17").Style = "Normal"

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.614.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E18").Style = "Normal"

# Row 19 - InternetComputer(DFINITY)
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.51%  "
$ws.Range("E19").Style = "Normal"

# Row 20 - ShibaInu
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000101"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.96%  "
$ws.Range("E20").Style = "Normal"

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +11.09%  "
$ws.Range("E21").Style = "Normal"

# Row 22 - Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.64%  "
$ws.Range("E22").Style = "Normal"

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.06"
$ws.Range("D23").Style = "Normal"

# Row 24 - PancakeSwap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.35%  "
$ws.Range("E24").Style = "Normal"

# Row 25 - ImmutableX
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +16.00%  "
$ws.Range("E25").Style = "Normal"

# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +35.04%  "
$ws.Range("E26").Style = "Normal"

# Row 27 - Dai
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E27").Style = "Normal"

# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.40%  "
$ws.Range("E28").Style = "Normal"

# Row 29 - InjectiveProtocol
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E29").Style = "Normal"

# Row 30 - Toncoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.50%  "
$ws.Range("E30").Style = "Normal"

# Row 31 - Filecoin
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +12.35%  "
$ws.Range("E31").Style = "Normal"

# Row 32 - LidoDAOToken
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E32").Style = "Normal"

# Row 33 - ARBITRUM
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +22.99%  "
$ws.Range("E33").Style = "Normal"

# Row 34 - WEMIXToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.24%  "
$ws.Range("E34").Style = "Normal"

# Row 35 - Hedera
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0835"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.41%  "
$ws.Range("E35").Style = "Normal"

# Row 36 - Monero
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("E36").Style = "Normal"

# Row 37 - Kaspa
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.46%  "
$ws.Range("E37").Style = "Normal"

# Row 38 - Stellar
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.91%  "
$ws.Range("E38").Style = "Normal"

# Row 39 - RenderToken->Celestia
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Celestia"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.34%  "
$ws.Range("E39").Style = "Normal"

# Row 40 - Celestia->RenderToken
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "RenderToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.53%  "
$ws.Range("E40").Style = "Normal"

# Row 41 - NEARProtocol
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +13.19%  "
$ws.Range("E41").Style = "Normal"

# Row 42 - VeChain
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0324"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.90%  "
$ws.Range("E42").Style = "Normal"

# Row 43 - Maker
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.054.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.51%  "
$ws.Range("E43").Style = "Normal"

# Row 44 - EnergySwap
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +31.43%  "
$ws.Range("E44").Style = "Normal"

# Row 45 - FirstDigitalUSD
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E45").Style = "Normal"

# Row 46 - BitcoinSV
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("E46").Style = "Normal"

# Row 47 - Stacks->FraxShare
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "FraxShare"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +10.56%  "
$ws.Range("E47").Style = "Normal"

# Row 48 - FraxShare->Stacks
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Stacks"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E48").Style = "Normal"

# Row 49 - Aave
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +12.34%  "
$ws.Range("E49").Style = "Normal"

# Row 50 - Algorand
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.17%  "
$ws.Range("E50").Style = "Normal"

# Row 51 - RocketPoolETH
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.866.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.43%  "
$ws.Range("E51").Style = "Normal"

